# The authored commit swaps the OOXML theme parts: the deck's main
# design theme (currently the "Integral" colour scheme, physically
# stored as ppt/theme/theme2.xml) is replaced by the stock "Office
# Theme" colour scheme that used to live in ppt/theme/theme1.xml (the
# notes master's theme, otherwise untouched by the rest of the deck).
#
# The font scheme and format scheme are identical between the two
# themes, so only the 10 colour-scheme slots that differ (dk2, lt2,
# accent1-6, hlink, folHlink -- dk1/lt1 are black/white in both) need
# to move across. We drive this through
# Slide.ThemeColorScheme.Colors(i).RGB, which is the PowerPoint OM's
# supported surface for editing the presentation's theme colours.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

function Set-ThemeColor {
    param(
        [int]$Index,
        [string]$HexRGB
    )
    $r = [Convert]::ToInt32($HexRGB.Substring(0,2), 16)
    $g = [Convert]::ToInt32($HexRGB.Substring(2,2), 16)
    $b = [Convert]::ToInt32($HexRGB.Substring(4,2), 16)
    # PowerPoint's ColorFormat.RGB is a Windows COLORREF: 0x00BBGGRR
    $val = ($b * 65536) + ($g * 256) + $r
    $cs.Colors($Index).RGB = $val
}

# Theme colour scheme slot order: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1
# 6=accent2 7=accent3 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
# dk1 (000000) and lt1 (FFFFFF) are unchanged between the two themes.
Set-ThemeColor 3  "44546A"   # dk2
Set-ThemeColor 4  "E7E6E6"   # lt2
Set-ThemeColor 5  "5B9BD5"   # accent1
Set-ThemeColor 6  "ED7D31"   # accent2
Set-ThemeColor 7  "A5A5A5"   # accent3
Set-ThemeColor 8  "FFC000"   # accent4
Set-ThemeColor 9  "4472C4"   # accent5
Set-ThemeColor 10 "70AD47"   # accent6
Set-ThemeColor 11 "0563C1"   # hlink
Set-ThemeColor 12 "954F72"   # folHlink
